$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 4
$ws.Range("I3").Value = 2.44
$ws.Range("T3").Value = 1.77
$ws.Range("U3").Value = 2.04
$ws.Range("W3").Value = 1.33
$ws.Range("AC3").Value = 9.6
$ws.Range("AO3").Value = 980
$ws.Range("F4").Value = 1.76
$ws.Range("G4").Value = 2.02
$ws.Range("J4").Value = 3.05
$ws.Range("Q4").Value = 2.42
$ws.Range("U4").Value = 1.55
$ws.Range("W4").Value = 1.98
$ws.Range("I5").Value = 2.96
$ws.Range("N5").Value = 6.4
$ws.Range("O5").Value = 1.14
$ws.Range("P5").Value = 2.82
$ws.Range("R5").Value = 1.74
$ws.Range("V5").Value = 1.51
$ws.Range("W5").Value = 1.58
$ws.Range("Y5").Value = 25
$ws.Range("AB5").Value = 24
$ws.Range("AF5").Value = 28
$ws.Range("AO5").Value = 15
$ws.Range("H6").Value = 16
$ws.Range("J6").Value = 6.8
$ws.Range("W6").Value = 4.7
$ws.Range("H7").Value = 1.33
$ws.Range("I7").Value = 1.48
$ws.Range("J7").Value = 4.2
$ws.Range("N7").Value = 3.25
$ws.Range("O7").Value = 1.26
$ws.Range("T7").Value = 2.3
$ws.Range("U7").Value = 1.61
$ws.Range("V7").Value = 3.05
$ws.Range("G8").Value = 3.6
$ws.Range("H8").Value = 2.32
$ws.Range("O8").Value = 1.37
$ws.Range("S8").Value = 3.85
$ws.Range("V8").Value = 1.63
$ws.Range("W8").Value = 1.38
$ws.Range("Q9").Value = 1.71
$ws.Range("F10").Value = 1.87
$ws.Range("G10").Value = 1.91
$ws.Range("I10").Value = 5.3
$ws.Range("J10").Value = 3.55
$ws.Range("K10").Value = 3.75
$ws.Range("N10").Value = 3.7
$ws.Range("O10").Value = 1.32
$ws.Range("P10").Value = 1.92
$ws.Range("Q10").Value = 1.94
$ws.Range("R10").Value = 1.34
$ws.Range("S10").Value = 3.35
$ws.Range("T10").Value = 1.8
$ws.Range("U10").Value = 2.08
$ws.Range("W10").Value = 2.1
$ws.Range("AC10").Value = 8.4
$ws.Range("AL10").Value = 980
$ws.Range("AM10").Value = 130
$ws.Range("G11").Value = 4.3
$ws.Range("I11").Value = 2.3
$ws.Range("J11").Value = 3.45
$ws.Range("N11").Value = 4
$ws.Range("V11").Value = 1.76
$ws.Range("W11").Value = 1.31
$ws.Range("AG11").Value = 18.5
$ws.Range("F12").Value = 1.79
$ws.Range("G12").Value = 1.92
$ws.Range("I12").Value = 5.7
$ws.Range("J12").Value = 3.65
$ws.Range("K12").Value = 4.3
$ws.Range("L12").Value = 1.36
$ws.Range("M12").Value = 1.08
$ws.Range("N12").Value = 3.3
$ws.Range("P12").Value = 1.79
$ws.Range("Q12").Value = 1.91
$ws.Range("R12").Value = 1.29
$ws.Range("S12").Value = 3.75
$ws.Range("T12").Value = 1.79
$ws.Range("U12").Value = 1.8
$ws.Range("V12").Value = 1.22
$ws.Range("W12").Value = 2.08
$ws.Range("X12").Value = 15.5
$ws.Range("AC12").Value = 10
$ws.Range("AN12").Value = 17.5
$ws.Range("AO12").Value = 130
$ws.Range("F13").Value = 1.31
$ws.Range("H13").Value = 10
$ws.Range("S13").Value = 2.04
$ws.Range("W13").Value = 3.95
$ws.Range("F14").Value = 2.86
$ws.Range("G14").Value = 3.4
$ws.Range("I14").Value = 2.52
$ws.Range("K14").Value = 4.9
$ws.Range("U14").Value = 2.4
$ws.Range("V14").Value = 1.65
$ws.Range("W14").Value = 1.41
$ws.Range("F16").Value = 3.3
$ws.Range("G16").Value = 5.1
$ws.Range("H16").Value = 1.43
$ws.Range("I16").Value = 2.76
$ws.Range("J16").Value = 2.54
$ws.Range("N16").Value = 1.25
$ws.Range("P16").Value = 1.24
$ws.Range("Q16").Value = 1.01
$ws.Range("S16").Value = 1.01
$ws.Range("V16").Value = 1.57
$ws.Range("W16").Value = 1.24
$ws.Range("F17").Value = 1.94
$ws.Range("G17").Value = 2.56
$ws.Range("H17").Value = 3.35
$ws.Range("I17").Value = 4.7
$ws.Range("J17").Value = 3.1
$ws.Range("K17").Value = 5.4
$ws.Range("Q17").Value = 1.98
$ws.Range("S17").Value = 1.99
$ws.Range("V17").Value = 1.27
$ws.Range("W17").Value = 1.64
$ws.Range("N18").Value = 3.45
$ws.Range("H20").Value = 4.3
$ws.Range("J20").Value = 3.4
$ws.Range("P21").Value = 2.02
$ws.Range("Q21").Value = 1.81
$ws.Range("AF21").Value = 65
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 2.14
$ws.Range("H22").Value = 3.8
$ws.Range("I22").Value = 4.4
$ws.Range("K22").Value = 4
$ws.Range("L22").Value = 1.32
$ws.Range("N22").Value = 3.8
$ws.Range("Q22").Value = 1.87
$ws.Range("T22").Value = 1.73
$ws.Range("U22").Value = 2.12
$ws.Range("V22").Value = 1.33
$ws.Range("W22").Value = 1.87
$ws.Range("AD22").Value = 20
$ws.Range("AF22").Value = 14
$ws.Range("AJ22").Value = 30
$ws.Range("AK22").Value = 23
$ws.Range("N23").Value = 1.59
$ws.Range("Q23").Value = 2.44
$ws.Range("F24").Value = 1.62
$ws.Range("G24").Value = 1.74
$ws.Range("H24").Value = 5.5
$ws.Range("J24").Value = 3.9
$ws.Range("K24").Value = 4.4
$ws.Range("R24").Value = 1.27
$ws.Range("W24").Value = 2.34
$ws.Range("X24").Value = 22
